# econ_activity ETL: append the next month (Feb/2021) row to the ICVA
# series and restore the on-screen selection to where the analyst left
# off after typing it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- copy the formatting (date / percentage styles) from the previous
#     last row down into the new row before writing values into it ----
$ws.Range("A98:E98").Copy() | Out-Null
$ws.Range("A99:E99").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- append the new month's figures (2021-02-01, serial 44228) -------
$ws.Range("A99").Value = 44228
$ws.Range("B99").Value = -0.099
$ws.Range("C99").Value = -0.072
$ws.Range("D99").Value = -0.171
$ws.Range("E99").Value = -0.121

# --- leave the selection where the user ended up after data entry ----
[void]$ws.Range("I95").Select()
